$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sprint1")
$ws2 = $wb.Worksheets.Item("Sprint2")

# Fill new rows 8 and 9 on Sprint2 with "N/A" across columns A:P
for ($col = 1; $col -le 16; $col++) {
    $ws2.Cells.Item(8, $col).Value = "N/A"
    $ws2.Cells.Item(9, $col).Value = "N/A"
}

# Update selection / active sheet state
$ws1.Range("C13").Select()
$ws2.Range("A8:P9").Select()
